# Weekly Fruta/Hortaliza update: a new (most recent) price observation is
# prepended to the existing data block. This pushes every existing data
# row (119-144) down by one row (120-145); row 119 is freed up for the
# new observation. Using Rows.Insert (Shift:=xlShiftDown semantics by
# default) reproduces that shift while Excel automatically carries the
# row-above formatting (including the custom date style on column D)
# onto the freshly inserted row, exactly like the diff shows (s="2" kept
# on D119/D145).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 119, pushing rows 119:144 down to 120:145.
$ws.Rows("119:119").Insert()

# Populate the newly freed row 119 with the latest weekly observation.
$ws.Range("A119").Value = 7
$ws.Range("B119").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C119").Value = "Ñuble"
$ws.Range("D119").Value = 44511
$ws.Range("E119").Value = 16
$ws.Range("F119").Value = 100112017
$ws.Range("G119").Value = "Apio"
$ws.Range("H119").Value = "Americana (o)"
$ws.Range("I119").Value = "Primera"
$ws.Range("J119").Value = 80
$ws.Range("K119").Value = 8000
$ws.Range("L119").Value = 9000
$ws.Range("M119").Value = 8500
$ws.Range("N119").Value = "`$/docena de matas"
$ws.Range("O119").Value = "Provincia del Elquí"
$ws.Range("P119").Value = 1417
$ws.Range("Q119").Value = 6
$ws.Range("R119").Value = "Hortaliza"
